# Power unit-converter workbook update
# - Insert two new rows (Default From Row / Default To Row) after the existing
#   "Image" row, pushing the Units table down.
# - Insert a new (empty-header) column C, pushing the old "Name comments" /
#   helper-text column out to D.
# - Populate the new column C with plain-language unit names for each unit row.
# - Adjust a couple of row heights to match the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing column B width before we insert a new column, so the
# freshly inserted column C can be given the same display width.
$colBWidth = $ws.Columns("B").ColumnWidth

# ---------------------------------------------------------------------------
# Structural edits: insert rows first (so row numbers below are the final,
# post-insert numbers), then insert the new column.
# ---------------------------------------------------------------------------

# Two new rows go in right after row 3 ("Image"), before the old row 4
# ("Units" header), which becomes row 6 once both rows are inserted.
$ws.Rows("4:5").Insert()

# A new column goes in at C, pushing the old helper/description column (old
# C) out to D.
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = $colBWidth

# ---------------------------------------------------------------------------
# New row 4 / row 5 content ("Default From Row" / "Default To Row").
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Default From Row"
$ws.Range("B4").Value = 8
$ws.Range("A5").Value = "Default To Row"
$ws.Range("B5").Value = 9

# ---------------------------------------------------------------------------
# New column C labels for each unit row (plain-language names).
# ---------------------------------------------------------------------------
$st = $wb.Styles.Add("Normal 2")
$st.Font.Name = "Calibri"
$st.Font.Size = 11

$ws.Range("C8:C12").Style = "Normal 2"
$ws.Range("C8").Value = "kilowatt"
$ws.Range("C9").Value = "horsepower"
$ws.Range("C10").Value = "watt"
$ws.Range("C11").Value = "foot-pound force per minute"
$ws.Range("C12").Value = "foot-pound force per second"

# ---------------------------------------------------------------------------
# Row-height tweaks to match the refreshed layout.
# ---------------------------------------------------------------------------
$ws.Rows("6").RowHeight = 16.5
$ws.Rows("9").RowHeight = 51
$ws.Rows("10:12").RowHeight = 15

# ---------------------------------------------------------------------------
# Selection, matching the saved view in the final workbook.
# ---------------------------------------------------------------------------
[void]$ws.Range("C5").Select()
